$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 333362660
$ws.Range("J69").Value = 333362660
$ws.Range("L69").Value = 1000087980
$ws.Range("N69").Value = -1000089728
$ws.Range("H70").Value = 5250.2075
$ws.Range("J70").Value = 5958.089
$ws.Range("L70").Value = 17874.267
$ws.Range("N70").Value = -18414.267
$ws.Range("H72").Value = 333362660
$ws.Range("J72").Value = 333362660
$ws.Range("L72").Value = 3000263940
$ws.Range("N72").Value = -3000272676
$ws.Range("H73").Value = 5250.2075
$ws.Range("J73").Value = 5958.089
$ws.Range("L73").Value = 17874.267
$ws.Range("N73").Value = -19746.267
$ws.Range("H76").Value = 6804.9375
$ws.Range("I76").Value = 5282
$ws.Range("K76").Value = 5282
$ws.Range("M76").Value = -4967
$ws.Range("H79").Value = 6804.9375
$ws.Range("I79").Value = 5282
$ws.Range("K79").Value = 5282
$ws.Range("M79").Value = -4190
$ws.Range("H80").Value = 682.7
$ws.Range("I80").Value = 434.85715
$ws.Range("J80").Value = 1261
$ws.Range("K80").Value = 1304.57145
$ws.Range("L80").Value = 3783
$ws.Range("M80").Value = -306.5714499999999
$ws.Range("N80").Value = -5779
$ws.Range("H83").Value = 682.7
$ws.Range("I83").Value = 434.85715
$ws.Range("J83").Value = 1261
$ws.Range("K83").Value = 3913.71435
$ws.Range("L83").Value = 11349
$ws.Range("M83").Value = 1078.28565
$ws.Range("N83").Value = -21333
$ws.Range("H135").Value = 550.45
$ws.Range("I135").Value = 578
$ws.Range("K135").Value = 5202
$ws.Range("M135").Value = -2667
$ws.Range("H138").Value = 3407.7856
$ws.Range("I138").Value = 1409.3889
$ws.Range("K138").Value = 4228.1667
$ws.Range("M138").Value = 911.8333000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5131634
$ws.Range("I45").Value = 7694213
$ws.Range("J45").Value = 6476.5
$ws.Range("K45").Value = 7694213
$ws.Range("L45").Value = 6476.5
$ws.Range("M45").Value = -7693836
$ws.Range("N45").Value = -7230.5
$ws.Range("H122").Value = 1606565.5
$ws.Range("I122").Value = 3809.8
$ws.Range("K122").Value = 11429.4
$ws.Range("M122").Value = -8979.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4234507
$ws.Range("I105").Value = 4536864.5
$ws.Range("K105").Value = 4536864.5
$ws.Range("M105").Value = -4535117.5
$ws.Range("H134").Value = 10442.647
$ws.Range("I134").Value = 10347.593
$ws.Range("K134").Value = 31042.779
$ws.Range("M134").Value = -28507.779

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4006.54
$ws.Range("I31").Value = 1614.5862
$ws.Range("J31").Value = 4983.535
$ws.Range("K31").Value = 1614.5862
$ws.Range("L31").Value = 4983.535
$ws.Range("M31").Value = -1319.5862
$ws.Range("N31").Value = -5573.535
$ws.Range("H34").Value = 4006.54
$ws.Range("I34").Value = 1614.5862
$ws.Range("J34").Value = 4983.535
$ws.Range("K34").Value = 1614.5862
$ws.Range("L34").Value = 4983.535
$ws.Range("M34").Value = -1412.5862
$ws.Range("N34").Value = -5387.535
$ws.Range("H122").Value = 2384.5
$ws.Range("J122").Value = 3187.6667
$ws.Range("L122").Value = 9563.000100000001
$ws.Range("N122").Value = -14463.0001
$ws.Range("H134").Value = 11354.16
$ws.Range("I134").Value = 10427.615
$ws.Range("J134").Value = 12357.917
$ws.Range("K134").Value = 31282.845
$ws.Range("L134").Value = 37073.751
$ws.Range("M134").Value = -28747.845
$ws.Range("N134").Value = -42143.751

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1202.7333
$ws.Range("I3").Value = 753.4167
$ws.Range("K3").Value = 2260.2501
$ws.Range("M3").Value = -2148.2501
$ws.Range("H11").Value = 7136.5386
$ws.Range("I11").Value = 2971.875
$ws.Range("J11").Value = 13800
$ws.Range("K11").Value = 8915.625
$ws.Range("L11").Value = 41400
$ws.Range("M11").Value = -8775.625
$ws.Range("N11").Value = -41680
$ws.Range("H107").Value = 1475.25
$ws.Range("I107").Value = 881.2
$ws.Range("K107").Value = 2643.6
$ws.Range("M107").Value = -723.6000000000004
$ws.Range("H131").Value = 26044580
$ws.Range("I131").Value = 41667268
$ws.Range("J131").Value = 23812768
$ws.Range("K131").Value = 125001804
$ws.Range("L131").Value = 71438304
$ws.Range("M131").Value = -124996764
$ws.Range("N131").Value = -71448384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1533584.5
$ws.Range("J80").Value = 4950.8335
$ws.Range("L80").Value = 4950.8335
$ws.Range("N80").Value = -6946.8335
$ws.Range("H83").Value = 1533584.5
$ws.Range("J83").Value = 4950.8335
$ws.Range("L83").Value = 24754.1675
$ws.Range("N83").Value = -34738.1675
$ws.Range("H102").Value = 4249826.5
$ws.Range("I102").Value = 6537313.5
$ws.Range("J102").Value = 1657341.4
$ws.Range("K102").Value = 6537313.5
$ws.Range("L102").Value = 1657341.4
$ws.Range("M102").Value = -6535691.5
$ws.Range("N102").Value = -1660585.4
$ws.Range("H129").Value = 39889.5
$ws.Range("J129").Value = 39889.5
$ws.Range("L129").Value = 39889.5
$ws.Range("N129").Value = -49889.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2239
$ws.Range("I16").Value = 2239
$ws.Range("K16").Value = 2239
$ws.Range("M16").Value = -2069
$ws.Range("H40").Value = 5446.6
$ws.Range("I40").Value = 4462.933
$ws.Range("J40").Value = 8397.6
$ws.Range("K40").Value = 4462.933
$ws.Range("L40").Value = 8397.6
$ws.Range("M40").Value = -4326.933
$ws.Range("N40").Value = -8669.6
$ws.Range("H68").Value = 4918.125
$ws.Range("I68").Value = 3529
$ws.Range("J68").Value = 7233.3335
$ws.Range("K68").Value = 3529
$ws.Range("L68").Value = 7233.3335
$ws.Range("M68").Value = -2780
$ws.Range("N68").Value = -8731.333500000001
$ws.Range("H71").Value = 4918.125
$ws.Range("I71").Value = 3529
$ws.Range("J71").Value = 7233.3335
$ws.Range("K71").Value = 17645
$ws.Range("L71").Value = 36166.6675
$ws.Range("M71").Value = -13901
$ws.Range("N71").Value = -43654.6675
$ws.Range("H100").Value = 3753.2693
$ws.Range("I100").Value = 3703.64
$ws.Range("K100").Value = 3703.64
$ws.Range("M100").Value = -3162.64
$ws.Range("H122").Value = 6946.25
$ws.Range("I122").Value = 5081.5713
$ws.Range("J122").Value = 9556.799999999999
$ws.Range("K122").Value = 15244.7139
$ws.Range("L122").Value = 28670.4
$ws.Range("M122").Value = -12794.7139
$ws.Range("N122").Value = -33570.39999999999
$ws.Range("H136").Value = 43798.69
$ws.Range("I136").Value = 52322.24
$ws.Range("K136").Value = 156966.72
$ws.Range("M136").Value = -154416.72

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10240
$ws.Range("I62").Value = 14427.667
$ws.Range("J62").Value = 9097.909
$ws.Range("K62").Value = 14427.667
$ws.Range("L62").Value = 9097.909
$ws.Range("M62").Value = -13803.667
$ws.Range("N62").Value = -10345.909
$ws.Range("H65").Value = 10240
$ws.Range("I65").Value = 14427.667
$ws.Range("J65").Value = 9097.909
$ws.Range("K65").Value = 72138.33499999999
$ws.Range("L65").Value = 45489.545
$ws.Range("M65").Value = -69018.33499999999
$ws.Range("N65").Value = -51729.545
$ws.Range("H81").Value = 15153251
$ws.Range("I81").Value = 18519584
$ws.Range("K81").Value = 37039168
$ws.Range("M81").Value = -37038107
$ws.Range("H84").Value = 15153251
$ws.Range("I84").Value = 18519584
$ws.Range("K84").Value = 185195840
$ws.Range("M84").Value = -185190536

Write-Host "Applied all cell updates"